# Update cryptocurrency price/volume data (GitHub Actions refresh, Fri Apr 26 15:40:18 UTC 2024).
# Row 31<->32 and 46<->47 also swap coin identity (re-ranked), matching the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '63.432.55'
$ws.Cells.Item(2, 5).Value = '  -0.85%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.108.85'
$ws.Cells.Item(3, 5).Value = '  -1.44%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.20%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '597.79'
$ws.Cells.Item(5, 5).Value = '  -1.88%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.82'
$ws.Cells.Item(6, 5).Value = '  -2.80%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.25%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.105.43'
$ws.Cells.Item(8, 5).Value = '  -1.49%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.518'
$ws.Cells.Item(9, 5).Value = '  -0.97%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -1.97%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.33'
$ws.Cells.Item(11, 5).Value = '  -0.99%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.45%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -1.23%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '35.04'
$ws.Cells.Item(14, 5).Value = '  -1.22%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.625.33'
$ws.Cells.Item(15, 5).Value = '  -0.74%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +2.39%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '63.661.92'
$ws.Cells.Item(17, 5).Value = '  -0.09%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.110.23'
$ws.Cells.Item(18, 5).Value = '  -0.78%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.77'
$ws.Cells.Item(19, 5).Value = '  -1.61%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '480.43'
$ws.Cells.Item(20, 5).Value = '  +0.80%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '14.55'
$ws.Cells.Item(21, 5).Value = '  -0.29%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -1.68%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.58'
$ws.Cells.Item(23, 5).Value = '  -4.86%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '86.96'
$ws.Cells.Item(24, 5).Value = '  +4.35%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '13.27'
$ws.Cells.Item(25, 5).Value = '  -3.36%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.01%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -2.91%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.24'
$ws.Cells.Item(28, 5).Value = '  -2.97%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.05'
$ws.Cells.Item(29, 5).Value = '  -1.07%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -2.73%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.00'
$ws.Cells.Item(31, 5).Value = '  +0.09%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '27.00'
$ws.Cells.Item(32, 5).Value = '  +2.70%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -8.67%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.32%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -2.69%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.30%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -4.10%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '52.44'
$ws.Cells.Item(38, 5).Value = '  -0.67%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.91'
$ws.Cells.Item(39, 5).Value = '  -2.80%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '435.62'
$ws.Cells.Item(40, 5).Value = '  -5.17%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0391'
$ws.Cells.Item(41, 5).Value = '  -1.02%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.49%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.90%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '2.856.14'
$ws.Cells.Item(44, 5).Value = '  -0.20%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -3.33%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Fetch.AI'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.20'
$ws.Cells.Item(46, 5).Value = '  -4.52%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ThetaToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.42'
$ws.Cells.Item(47, 5).Value = '  +0.68%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.07%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '25.70'
$ws.Cells.Item(49, 5).Value = '  -2.77%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.41%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '121.26'
$ws.Cells.Item(51, 5).Value = '  +2.03%  '
